# Update crypto price/volume data per Wed May  1 11:43:54 UTC 2024 GitHub Actions run
$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Sheet1")

# Reference style of an unstyled data cell (no border/bold), used to strip the
# "quote prefix" style Excel applies when a value is entered as explicit text,
# so modified cells keep matching the (un-styled) look of the rest of the sheet.
$plainStyle = $ws.Range("B2").Style

$ws.Range('D2').Value = "'57.568.35"
$ws.Range('D2').Style = $plainStyle
$ws.Range('E2').Value = "'  -5.93%  "
$ws.Range('E2').Style = $plainStyle
$ws.Range('D3').Value = "'2.885.36"
$ws.Range('D3').Style = $plainStyle
$ws.Range('E3').Value = "'  -4.19%  "
$ws.Range('E3').Style = $plainStyle
$ws.Range('E4').Value = "'  -0.11%  "
$ws.Range('E4').Style = $plainStyle
$ws.Range('D5').Value = "'548.17"
$ws.Range('D5').Style = $plainStyle
$ws.Range('E5').Value = "'  -4.07%  "
$ws.Range('E5').Style = $plainStyle
$ws.Range('D6').Value = "'121.72"
$ws.Range('D6').Style = $plainStyle
$ws.Range('E6').Value = "'  -5.51%  "
$ws.Range('E6').Style = $plainStyle
$ws.Range('E7').Value = "'  +0.09%  "
$ws.Range('E7').Style = $plainStyle
$ws.Range('D8').Value = "'2.882.32"
$ws.Range('D8').Style = $plainStyle
$ws.Range('E8').Value = "'  -4.22%  "
$ws.Range('E8').Style = $plainStyle
$ws.Range('D9').Value = "'0.493"
$ws.Range('D9').Style = $plainStyle
$ws.Range('E9').Value = "'  -0.78%  "
$ws.Range('E9').Style = $plainStyle
$ws.Range('E10').Value = "'  -8.18%  "
$ws.Range('E10').Style = $plainStyle
$ws.Range('E11').Value = "'  -9.06%  "
$ws.Range('E11').Style = $plainStyle
$ws.Range('E12').Value = "'  +0.14%  "
$ws.Range('E12').Style = $plainStyle
$ws.Range('E13').Value = "'  -7.60%  "
$ws.Range('E13').Style = $plainStyle
$ws.Range('D14').Value = "'31.60"
$ws.Range('D14').Style = $plainStyle
$ws.Range('E14').Value = "'  -4.94%  "
$ws.Range('E14').Style = $plainStyle
$ws.Range('E15').Value = "'  -0.47%  "
$ws.Range('E15').Style = $plainStyle
$ws.Range('D16').Value = "'3.364.97"
$ws.Range('D16').Style = $plainStyle
$ws.Range('E16').Value = "'  -4.14%  "
$ws.Range('E16').Style = $plainStyle
$ws.Range('D17').Value = "'2.887.88"
$ws.Range('D17').Style = $plainStyle
$ws.Range('E17').Value = "'  -4.43%  "
$ws.Range('E17').Style = $plainStyle
$ws.Range('B18').Value = "'WrappedBTC"
$ws.Range('B18').Style = $plainStyle
$ws.Range('C18').Value = "'https://coinranking.com/coin/x4WXHge-vvFY+wrappedbtc-wbtc"
$ws.Range('C18').Style = $plainStyle
$ws.Range('D18').Value = "'57.484.59"
$ws.Range('D18').Style = $plainStyle
$ws.Range('E18').Value = "'  -6.27%  "
$ws.Range('E18').Style = $plainStyle
$ws.Range('B19').Value = "'Polkadot"
$ws.Range('B19').Style = $plainStyle
$ws.Range('C19').Value = "'https://coinranking.com/coin/25W7FG7om+polkadot-dot"
$ws.Range('C19').Style = $plainStyle
$ws.Range('D19').Value = "'6.46"
$ws.Range('D19').Style = $plainStyle
$ws.Range('E19').Value = "'  +3.01%  "
$ws.Range('E19').Style = $plainStyle
$ws.Range('D20').Value = "'406.77"
$ws.Range('D20').Style = $plainStyle
$ws.Range('E20').Value = "'  -7.11%  "
$ws.Range('E20').Style = $plainStyle
$ws.Range('D21').Value = "'12.81"
$ws.Range('D21').Style = $plainStyle
$ws.Range('E21').Value = "'  -3.26%  "
$ws.Range('E21').Style = $plainStyle
$ws.Range('E22').Value = "'  -1.46%  "
$ws.Range('E22').Style = $plainStyle
$ws.Range('D23').Value = "'6.74"
$ws.Range('D23').Style = $plainStyle
$ws.Range('E23').Value = "'  -6.39%  "
$ws.Range('E23').Style = $plainStyle
$ws.Range('D24').Value = "'12.56"
$ws.Range('D24').Style = $plainStyle
$ws.Range('E24').Value = "'  -0.57%  "
$ws.Range('E24').Style = $plainStyle
$ws.Range('D25').Value = "'76.81"
$ws.Range('D25').Style = $plainStyle
$ws.Range('E25').Value = "'  -3.78%  "
$ws.Range('E25').Style = $plainStyle
$ws.Range('E26').Value = "'  +0.15%  "
$ws.Range('E26').Style = $plainStyle
$ws.Range('D27').Value = "'0.998"
$ws.Range('D27').Style = $plainStyle
$ws.Range('E27').Value = "'  -0.23%  "
$ws.Range('E27').Style = $plainStyle
$ws.Range('E28').Value = "'  -2.63%  "
$ws.Range('E28').Style = $plainStyle
$ws.Range('D29').Value = "'7.18"
$ws.Range('D29').Style = $plainStyle
$ws.Range('E29').Value = "'  -2.23%  "
$ws.Range('E29').Style = $plainStyle
$ws.Range('E30').Value = "'  -3.53%  "
$ws.Range('E30').Style = $plainStyle
$ws.Range('D31').Value = "'6.02"
$ws.Range('D31').Style = $plainStyle
$ws.Range('E31').Value = "'  -4.00%  "
$ws.Range('E31').Style = $plainStyle
$ws.Range('D32').Value = "'24.57"
$ws.Range('D32').Style = $plainStyle
$ws.Range('E32').Value = "'  -3.75%  "
$ws.Range('E32').Style = $plainStyle
$ws.Range('D33').Value = "'0.0948"
$ws.Range('D33').Style = $plainStyle
$ws.Range('E33').Value = "'  -0.26%  "
$ws.Range('E33').Style = $plainStyle
$ws.Range('E34').Value = "'  -11.67%  "
$ws.Range('E34').Style = $plainStyle
$ws.Range('D35').Value = "'0.895"
$ws.Range('D35').Style = $plainStyle
$ws.Range('E35').Value = "'  -6.59%  "
$ws.Range('E35').Style = $plainStyle
$ws.Range('E36').Value = "'  -4.88%  "
$ws.Range('E36').Style = $plainStyle
$ws.Range('D37').Value = "'48.31"
$ws.Range('D37').Style = $plainStyle
$ws.Range('E37').Value = "'  -3.66%  "
$ws.Range('E37').Style = $plainStyle
$ws.Range('D38').Value = "'8.41"
$ws.Range('D38').Style = $plainStyle
$ws.Range('E38').Value = "'  +7.95%  "
$ws.Range('E38').Style = $plainStyle
$ws.Range('D39').Value = "'0.0₃0616"
$ws.Range('D39').Style = $plainStyle
$ws.Range('E39').Value = "'  -9.56%  "
$ws.Range('E39').Style = $plainStyle
$ws.Range('D40').Value = "'0.0342"
$ws.Range('D40').Style = $plainStyle
$ws.Range('E40').Value = "'  -6.95%  "
$ws.Range('E40').Style = $plainStyle
$ws.Range('E41').Value = "'  -3.91%  "
$ws.Range('E41').Style = $plainStyle
$ws.Range('D42').Value = "'2.606.26"
$ws.Range('D42').Style = $plainStyle
$ws.Range('E42').Value = "'  -1.88%  "
$ws.Range('E42').Style = $plainStyle
$ws.Range('D43').Value = "'355.63"
$ws.Range('D43').Style = $plainStyle
$ws.Range('E43').Value = "'  -5.07%  "
$ws.Range('E43').Style = $plainStyle
$ws.Range('E44').Value = "'  -5.71%  "
$ws.Range('E44').Style = $plainStyle
$ws.Range('D45').Value = "'0.999"
$ws.Range('D45').Style = $plainStyle
$ws.Range('E45').Value = "'  -0.01%  "
$ws.Range('E45').Style = $plainStyle
$ws.Range('D46').Value = "'117.81"
$ws.Range('D46').Style = $plainStyle
$ws.Range('E46').Value = "'  -3.36%  "
$ws.Range('E46').Style = $plainStyle
$ws.Range('D47').Value = "'0.228"
$ws.Range('D47').Style = $plainStyle
$ws.Range('E47').Value = "'  -3.54%  "
$ws.Range('E47').Style = $plainStyle
$ws.Range('E48').Value = "'  -0.77%  "
$ws.Range('E48').Style = $plainStyle
$ws.Range('D49').Value = "'1.92"
$ws.Range('D49').Style = $plainStyle
$ws.Range('E49').Value = "'  -2.21%  "
$ws.Range('E49').Style = $plainStyle
$ws.Range('D50').Value = "'22.44"
$ws.Range('D50').Style = $plainStyle
$ws.Range('E50').Value = "'  -4.44%  "
$ws.Range('E50').Style = $plainStyle
$ws.Range('D51').Value = "'0.127"
$ws.Range('D51').Style = $plainStyle
$ws.Range('E51').Value = "'  -2.25%  "
$ws.Range('E51').Style = $plainStyle
